$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a Number by
# Excel's input parser (plain decimals like '254.58') are written with
# the cell pre-formatted as Text ('@') so they round-trip as strings,
# matching the workbook's original inlineStr storage for these columns.

$ws.Range('D2').Value = '98.247.04'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '3.405.73'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '254.58'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '668.06'
$ws.Range('E6').Value = '  -2.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.51'
$ws.Range('E7').Value = '  +4.96%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.432'
$ws.Range('E8').Value = '  +0.72%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.06'
$ws.Range('E9').Value = '  +1.36%  '
$ws.Range('E10').Value = '  -0.05%  '
$ws.Range('D11').Value = '3.402.03'
$ws.Range('E11').Value = '  +0.35%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '45.70'
$ws.Range('E12').Value = '  +10.26%  '
$ws.Range('E13').Value = '  -2.16%  '
$ws.Range('D14').Value = '98.097.32'
$ws.Range('E14').Value = '  +0.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.15'
$ws.Range('E15').Value = '  -1.63%  '
$ws.Range('E16').Value = '  -1.67%  '
$ws.Range('D17').Value = '4.037.68'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '9.14'
$ws.Range('E18').Value = '  +2.91%  '
$ws.Range('D19').Value = '3.437.40'
$ws.Range('E19').Value = '  +0.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.44'
$ws.Range('E20').Value = '  +6.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.536'
$ws.Range('E21').Value = '  -5.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.43'
$ws.Range('E22').Value = '  +4.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '512.72'
$ws.Range('E23').Value = '  +1.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.43'
$ws.Range('E24').Value = '  +0.54%  '
$ws.Range('E25').Value = '  -0.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.86'
$ws.Range('E26').Value = '  +4.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '97.68'
$ws.Range('E27').Value = '  -2.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.48'
$ws.Range('E28').Value = '  -1.16%  '
$ws.Range('D29').Value = '3.586.79'
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '12.33'
$ws.Range('E30').Value = '  +7.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.90'
$ws.Range('E31').Value = '  +10.60%  '
$ws.Range('E32').Value = '  -3.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.189'
$ws.Range('E34').Value = '  -2.88%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.569'
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '29.20'
$ws.Range('E37').Value = '  -1.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.99'
$ws.Range('E38').Value = '  +0.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.50'
$ws.Range('E39').Value = '  -0.60%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '529.86'
$ws.Range('E40').Value = '  +1.08%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.154'
$ws.Range('E41').Value = '  +1.03%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.866'
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('B44').Value = 'WhiteBITCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '24.42'
$ws.Range('E44').Value = '  -1.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.76'
$ws.Range('E45').Value = '  +2.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0428'
$ws.Range('E46').Value = '  -1.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.68'
$ws.Range('E47').Value = '  -2.80%  '
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.69'
$ws.Range('E48').Value = '  -2.50%  '
$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.66'
$ws.Range('E49').Value = '  -1.08%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.24'
$ws.Range('E50').Value = '  +6.46%  '
$ws.Range('B51').Value = 'OKB'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '56.02'
$ws.Range('E51').Value = '  +0.29%  '
